$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 2056.5
$ws.Range("I42").Value = 2675.3333
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 8025.999899999999
$ws.Range("L42").Value = 600
$ws.Range("M42").Value = -7795.999899999999
$ws.Range("N42").Value = -1060

# Row 70
$ws.Range("H70").Value = 3826.1155
$ws.Range("I70").Value = 5363.5
$ws.Range("J70").Value = 2032.5
$ws.Range("K70").Value = 16090.5
$ws.Range("L70").Value = 6097.5
$ws.Range("M70").Value = -15820.5
$ws.Range("N70").Value = -6637.5

# Row 73
$ws.Range("H73").Value = 3826.1155
$ws.Range("I73").Value = 5363.5
$ws.Range("J73").Value = 2032.5
$ws.Range("K73").Value = 16090.5
$ws.Range("L73").Value = 6097.5
$ws.Range("M73").Value = -15154.5
$ws.Range("N73").Value = -7969.5

# Row 76
$ws.Range("H76").Value = 8697.154
$ws.Range("I76").Value = 17868
$ws.Range("J76").Value = 5318.421
$ws.Range("K76").Value = 17868
$ws.Range("L76").Value = 5318.421
$ws.Range("M76").Value = -17553
$ws.Range("N76").Value = -5948.421

# Row 79
$ws.Range("H79").Value = 8697.154
$ws.Range("I79").Value = 17868
$ws.Range("J79").Value = 5318.421
$ws.Range("K79").Value = 17868
$ws.Range("L79").Value = 5318.421
$ws.Range("M79").Value = -16776
$ws.Range("N79").Value = -7502.421

# Row 93
$ws.Range("H93").Value = 36323.332
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 36323.332
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 36323.332
$ws.Range("N93").Value = -41315.332

# Row 98
$ws.Range("H98").Value = 1536.2667
$ws.Range("I98").Value = 1541.8462
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 1541.8462
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -43.84619999999995
$ws.Range("N98").Value = -4496

# Row 111
$ws.Range("H111").Value = 1765.8572
$ws.Range("I111").Value = 1764.5
$ws.Range("J111").Value = 1766.4
$ws.Range("K111").Value = 5293.5
$ws.Range("L111").Value = 5299.200000000001
$ws.Range("M111").Value = -2226.5
$ws.Range("N111").Value = -11433.2

# Row 122
$ws.Range("H122").Value = 1536.2667
$ws.Range("I122").Value = 1541.8462
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4625.5386
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2175.5386
$ws.Range("N122").Value = -9400

# Row 125
$ws.Range("H125").Value = 558.4761999999999
$ws.Range("I125").Value = 645.7143
$ws.Range("J125").Value = 384
$ws.Range("K125").Value = 5811.428699999999
$ws.Range("L125").Value = 3456
$ws.Range("M125").Value = -3351.428699999999
$ws.Range("N125").Value = -8376

# Row 129
$ws.Range("H129").Value = 1027.9375
$ws.Range("I129").Value = 357.5
$ws.Range("J129").Value = 1251.4166
$ws.Range("K129").Value = 1072.5
$ws.Range("L129").Value = 3754.2498
$ws.Range("M129").Value = 3927.5
$ws.Range("N129").Value = -13754.2498

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = ""

# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = ""

# Row 45
$ws.Range("H45").Value = 1634.4445
$ws.Range("I45").Value = 1552.5
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 1552.5
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -1175.5
$ws.Range("N45").Value = -2454

# Row 122
$ws.Range("H122").Value = 2819.5417
$ws.Range("I122").Value = 3030.0588
$ws.Range("J122").Value = 2308.2856
$ws.Range("K122").Value = 9090.1764
$ws.Range("L122").Value = 6924.8568
$ws.Range("M122").Value = -6640.1764
$ws.Range("N122").Value = -11824.8568

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -360

# Row 134
$ws.Range("H134").Value = 32309.828
$ws.Range("I134").Value = 52160.2
$ws.Range("J134").Value = 5842.6665
$ws.Range("K134").Value = 156480.6
$ws.Range("L134").Value = 17527.9995
$ws.Range("M134").Value = -153945.6

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1776
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1776
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 5328
$ws.Range("N122").Value = -10228

$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1327
$ws.Range("N16").Value = ""

# Row 23
$ws.Range("H23").Value = 521.48
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 561.1739
$ws.Range("K23").Value = 195
$ws.Range("L23").Value = 1683.5217
$ws.Range("M23").Value = 40
$ws.Range("N23").Value = -2153.5217

# Row 38
$ws.Range("H38").Value = 43.785713
$ws.Range("I38").Value = 40.18182
$ws.Range("J38").Value = 57
$ws.Range("K38").Value = 120.54546
$ws.Range("L38").Value = 171
$ws.Range("M38").Value = 226.45454

# Row 113
$ws.Range("H113").Value = 521.35297
$ws.Range("I113").Value = 488.66666
$ws.Range("J113").Value = 599.8
$ws.Range("K113").Value = 1465.99998
$ws.Range("L113").Value = 1799.4
$ws.Range("M113").Value = 704.0000199999999
$ws.Range("N113").Value = -6139.4

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 800644
$ws.Range("I14").Value = 800644
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 800644
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -800476
$ws.Range("N14").Value = ""

# Row 20
$ws.Range("H20").Value = 9800
$ws.Range("I20").Value = 9800
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 9800
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -9555
$ws.Range("N20").Value = ""

# Row 58
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = ""

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3440
$ws.Range("I122").Value = 3350
$ws.Range("J122").Value = 3575
$ws.Range("K122").Value = 10050
$ws.Range("L122").Value = 10725
$ws.Range("M122").Value = -7600

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 7500
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -4760
$ws.Range("N20").Value = -10480

# Row 107
$ws.Range("H107").Value = 573.5
$ws.Range("I107").Value = 520
$ws.Range("J107").Value = 680.5
$ws.Range("K107").Value = 1560
$ws.Range("L107").Value = 2041.5
$ws.Range("M107").Value = 360
$ws.Range("N107").Value = -5881.5
